# The sheet holds monthly index values for years 2014-2017, one row per
# month, 12 rows per year, starting at row 2 (row 1 is the header).
# For every year-block of 12 rows, the last 3 rows (Oct, Nov, Dec) need to
# move to the front of that block, pushing Jan..Sep down by 3 rows -
# i.e. a rotation of each 12-row block by 3.
#
# Example for the first block (rows 2-13, originally months 01..12):
#   before: 01 02 03 04 05 06 07 08 09 10 11 12
#   after : 10 11 12 01 02 03 04 05 06 07 08 09

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$rowsPerYear = 12
$yearCount = 4
$totalRows = $rowsPerYear * $yearCount
$lastDataRow = $firstDataRow + $totalRows - 1
$colCount = 6

# Read the full data block (A2:F49) in one shot.
$srcRange = $ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastDataRow, $colCount))
$src = $srcRange.Value2

# Build the destination array with each 12-row year-block rotated so the
# last 3 rows come first.
$dst = New-Object 'object[,]' $totalRows,$colCount

for ($year = 0; $year -lt $yearCount; $year++) {
    $blockStart = $year * $rowsPerYear   # 0-based offset of this year's block

    for ($i = 0; $i -lt $rowsPerYear; $i++) {
        # source row within the block (0-based): rotate by 3
        $srcRowInBlock = ($i + ($rowsPerYear - 3)) % $rowsPerYear
        $srcRow1 = $blockStart + $srcRowInBlock + 1   # 1-based index into $src
        $dstRow0 = $blockStart + $i                   # 0-based index into $dst

        for ($c = 1; $c -le $colCount; $c++) {
            $dst[$dstRow0, $c - 1] = $src[$srcRow1, $c]
        }
    }
}

$dstRange = $ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastDataRow, $colCount))
$dstRange.Value = $dst

Write-Output "rotated $totalRows rows across $yearCount year-blocks"
